$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (std) to hold the "median" stat.
$ws.Columns.Item(4).Insert()

# Header for new column
$ws.Cells.Item(1, 4).Value2 = "median"

# Median values per rater row (rows 2..24)
$medians = @(91, 91, 91, 90, 90, 90, 90, 89, 89, 89, 89, 89, 88, 89, 89, 88, 89, 88, 88, 88, 87, 87, 87)

for ($i = 0; $i -lt $medians.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value2 = $medians[$i]
}

Write-Host "done"
